$d = $word.ActiveDocument

# 1) Update the first paragraph: append trailing spaces after "git init"
$p1 = $d.Paragraphs.Item(1)
$p1.Range.Text = "PS C:\Users\winso\Desktop\Tradata Website> git init                               "

# 2) Replace paragraphs 3..67 (everything between the "Reinitialized..." line
#    and the final trailing "PS ...>" prompt line) with the new transcript
#    content, using InsertXML so we can control exact run-level markup
#    (xml:space, true empty paragraphs, and the lastRenderedPageBreak marker).
$pStart = $d.Paragraphs.Item(3)
$pEnd = $d.Paragraphs.Item(67)
$r = $d.Range($pStart.Range.Start, $pEnd.Range.End)

$xml = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">PS C:\Users\winso\Desktop\Tradata Website&gt; git remote add origin </w:t></w:r></w:p><w:p><w:r><w:t>usage: git remote add [&lt;options&gt;] &lt;name&gt; &lt;url&gt;</w:t></w:r></w:p><w:p/><w:p><w:r><w:t xml:space="preserve">    -f, --[no-]fetch      fetch the remote branches</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    --[no-]tags           import all tags and associated objects when fetching</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">                          or do not fetch any tag at all (--no-tags)</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    -t, --[no-]track &lt;branch&gt;</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">                          branch(es) to track</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    -m, --[no-]master &lt;branch&gt;</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">                          master branch</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    --[no-]mirror[=(push|fetch)]</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">                          set up remote as a mirror to push to or fetch from</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>PS C:\Users\winso\Desktop\Tradata Website&gt; git remote -v</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">origin  https://github.com/Dannywinson1/tradata.git (fetch)</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">origin  https://github.com/Dannywinson1/tradata.git (push)</w:t></w:r></w:p><w:p><w:r><w:t>PS C:\Users\winso\Desktop\Tradata Website&gt; git add .</w:t></w:r></w:p><w:p><w:r><w:t>PS C:\Users\winso\Desktop\Tradata Website&gt; git commit -m "Initial commit - upload Tradata website"</w:t></w:r></w:p><w:p><w:r><w:t>[main 0c19ec2] Initial commit - upload Tradata website</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve"> 4 files changed, 8 insertions(+), 8 deletions(-)</w:t></w:r></w:p><w:p><w:r><w:t>PS C:\Users\winso\Desktop\Tradata Website&gt; git push -u origin main</w:t></w:r></w:p><w:p><w:r><w:t>Enumerating objects: 11, done.</w:t></w:r></w:p><w:p><w:r><w:t>Counting objects: 100% (11/11), done.</w:t></w:r></w:p><w:p><w:r><w:t>Delta compression using up to 12 threads</w:t></w:r></w:p><w:p><w:r><w:t>Compressing objects: 100% (6/6), done.</w:t></w:r></w:p><w:p><w:r><w:lastRenderedPageBreak/><w:t>Writing objects: 100% (6/6), 924 bytes | 924.00 KiB/s, done.</w:t></w:r></w:p><w:p><w:r><w:t>Total 6 (delta 5), reused 0 (delta 0), pack-reused 0 (from 0)</w:t></w:r></w:p><w:p><w:r><w:t>remote: Resolving deltas: 100% (5/5), completed with 5 local objects.</w:t></w:r></w:p><w:p><w:r><w:t>To https://github.com/Dannywinson1/tradata.git</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">   983cfb5..0c19ec2  main -&gt; main</w:t></w:r></w:p><w:p><w:r><w:t>branch ''main'' set up to track ''origin/main''.</w:t></w:r></w:p><w:p><w:r><w:t>PS C:\Users\winso\Desktop\Tradata Website&gt; git push origin main</w:t></w:r></w:p><w:p><w:r><w:t>Everything up-to-date</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$r.InsertXML($xml)
